# Weekly update: insert a new observation row at row 247 (pushing the
# existing rows 247-268 down to 248-269), mirroring the "Fruta / hortaliza,
# semanal" weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 247; this shifts rows
# 247..268 down to 248..269 and grows the sheet from 268 to 269 data rows.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with this week's new record.
$ws.Range("A247").Value = 8
$ws.Range("B247").Value = "Terminal La Palmera de La Serena"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44578
$ws.Range("E247").Value = 4
$ws.Range("F247").Value = 100114013
$ws.Range("G247").Value = "Zanahoria"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 600
$ws.Range("K247").Value = 6000
$ws.Range("L247").Value = 6500
$ws.Range("M247").Value = 6250
$ws.Range("N247").Value = "$/saco 20 kilos"
$ws.Range("O247").Value = "Provincia del Elquí"
$ws.Range("P247").Value = 312
$ws.Range("Q247").Value = 20
$ws.Range("R247").Value = "Hortaliza"
